# Updates cryptocurrency price/volume data cells (columns D and E) to
# match the latest scrape. Values are assigned with a leading apostrophe
# to force text interpretation (matching the source inlineStr type),
# then the cell style is reset to "Normal" so no stray number-format /
# quote-prefix style is left behind on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'58.497.35"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -1.26%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'2.483.46"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -1.65%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.02%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'526.00"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -2.15%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'133.79"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -3.28%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E8").Value = "'  -0.85%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'  -1.38%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'  -1.82%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'5.39"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +0.28%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = "'  -1.20%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'2.924.35"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -1.74%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'58.387.26"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -1.27%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'22.41"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -3.65%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("E16").Value = "'  -2.15%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'2.487.84"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Value = "'10.93"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -1.65%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'4.21"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -2.06%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'321.30"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -1.59%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = "'  +0.10%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'5.83"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -1.31%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'64.39"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -2.06%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'  -2.55%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'1.00"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Value = "'0.161"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -2.80%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'  -2.66%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'0.0₃0752"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -3.37%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'6.44"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -4.95%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = "'  -3.87%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'165.50"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -2.45%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = "'  -4.18%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = "'  -0.02%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = "'  -0.08%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'18.26"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -1.54%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = "'  -8.91%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'3.99"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -3.35%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = "'  -4.04%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.799"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -3.64%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = "'  -3.37%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'276.41"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -2.89%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'4.97"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -5.36%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.595"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -1.82%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'127.37"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -2.38%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("E45").Value = "'  -1.83%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'  -3.18%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'  -2.79%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'17.21"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -1.83%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'1.740.23"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -1.33%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.974"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -1.54%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "'  -1.75%  "
$ws.Range("E51").Style = "Normal"
